# Append 4 new log rows (time/cost) to the bottom of the sheet, mirroring
# the existing "time" / "cost" table layout (Github Auto Build commit at
# 2023-12-12 15:15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current data (row 238 -> 239).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

$newRows = @(
    @("2023-12-12 15:14:40", 0.0008),
    @("2023-12-12 15:14:56", 0.0006000000000000001),
    @("2023-12-12 15:15:10", 0.0006000000000000001),
    @("2023-12-12 15:15:22", 0.0004)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
